# Apply the edit described by the diff:
#  - Insert a new "A005" row into the node (sheet2) on-us group
#  - Insert two new "C004"/"C005" rows into the node (sheet2) off-us group
#  - Make the "node" sheet (sheet2) the active sheet / active tab
#  - Update the selection on sheet2 to F12

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Insert "A005 / on-us / orange" right after the existing "A004" row (old row 5 -> new row 6)
$ws2.Rows.Item(6).Insert()
$ws2.Cells.Item(6, 1).Value = "A005"
$ws2.Cells.Item(6, 2).Value = "on-us"
$ws2.Cells.Item(6, 3).Value = "orange"

# Insert "C004" and "C005" / off-us / blue rows before the trailing "V001" row
$ws2.Rows.Item(11).Insert()
$ws2.Rows.Item(11).Insert()
$ws2.Cells.Item(11, 1).Value = "C004"
$ws2.Cells.Item(11, 2).Value = "off-us"
$ws2.Cells.Item(11, 3).Value = "blue"
$ws2.Cells.Item(12, 1).Value = "C005"
$ws2.Cells.Item(12, 2).Value = "off-us"
$ws2.Cells.Item(12, 3).Value = "blue"

# Make "node" (sheet2) the active sheet/tab and set its selection
$ws2.Activate()
$ws2.Range("F12").Select()
